# Implement DDT (data-driven testing) support: rename the existing sheet
# to "LoginTest" and add a new "RegisterTest" sheet populated with sample
# registration data.

$wb = $excel.ActiveWorkbook

# Rename the first (and currently only) sheet.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "LoginTest"

# Add the new sheet right after LoginTest and rename it.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "RegisterTest"

# Header row.
$headers = @("firstname", "lastname", "telephone", "password", "newsletter_yes_or_no", "select_user_agreement")
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws2.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Sample registration rows.
$data = @(
    @("John",   "Doe",     1234567890, 12345, "no", "select"),
    @("Karter",  "Scott",   2345678910, 23456, "no", "select"),
    @("Jane",    "Miligan", 3456789101, 34567, "no", "select"),
    @("Karen",   "Moris",   4567891012, 45678, "no", "select"),
    @("Tod",     "Haris",   5678901234, 56789, "no", "select")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws2.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Match column widths to the best-fit size Excel would use for this content
# (first name, telephone, password, newsletter flag and agreement columns
# are all wider than the default and get auto-sized; the last-name column
# keeps the sheet's default width).
$ws2.Columns.Item(1).AutoFit() | Out-Null
$ws2.Columns.Item(3).AutoFit() | Out-Null
$ws2.Columns.Item(4).AutoFit() | Out-Null
$ws2.Columns.Item(5).AutoFit() | Out-Null
$ws2.Columns.Item(6).AutoFit() | Out-Null

# Page setup for the new sheet.
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Put the active selection on F1 of the new (now active/front) sheet.
$ws2.Range("F1").Select() | Out-Null
